# Auto-generated edit script: updates the cryptos list values
# (prices in column D and 1h volume % in column E, plus a couple of
# coin name/link swaps in column B/C) to match the target snapshot.
$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

$ws.Range("D2").Value = '43.786.21'
$ws.Range("E2").Value = '  +5.06%  '
$ws.Range("D3").Value = '2.279.88'
$ws.Range("E3").Value = '  +3.20%  '
$ws.Range("E4").Value = '  +0.08%  '
$ws.Range("D5").NumberFormat = "@"
$ws.Range("D5").Value = '233.63'
$ws.Range("E5").Value = '  +1.62%  '
$ws.Range("D6").NumberFormat = "@"
$ws.Range("D6").Value = '0.640'
$ws.Range("E6").Value = '  +3.59%  '
$ws.Range("D7").NumberFormat = "@"
$ws.Range("D7").Value = '65.67'
$ws.Range("E7").Value = '  +8.78%  '
$ws.Range("E8").Value = '  +0.00%  '
$ws.Range("D9").NumberFormat = "@"
$ws.Range("D9").Value = '0.428'
$ws.Range("E9").Value = '  +6.27%  '
$ws.Range("D10").NumberFormat = "@"
$ws.Range("D10").Value = '0.103'
$ws.Range("E10").Value = '  +16.48%  '
$ws.Range("D11").NumberFormat = "@"
$ws.Range("D11").Value = '57.66'
$ws.Range("E11").Value = '  +0.69%  '
$ws.Range("E12").Value = '  +17.79%  '
$ws.Range("E13").Value = '  +0.63%  '
$ws.Range("D14").Value = '2.616.95'
$ws.Range("E14").Value = '  +3.16%  '
$ws.Range("D15").NumberFormat = "@"
$ws.Range("D15").Value = '15.89'
$ws.Range("E15").Value = '  +3.16%  '
$ws.Range("D16").NumberFormat = "@"
$ws.Range("D16").Value = '5.96'
$ws.Range("E16").Value = '  +5.12%  '
$ws.Range("E17").Value = '  +4.94%  '
$ws.Range("D18").Value = '2.283.82'
$ws.Range("E18").Value = '  +3.00%  '
$ws.Range("D19").Value = '43.707.04'
$ws.Range("E19").Value = '  +4.93%  '
$ws.Range("D20").Value = '0.0₃0994'
$ws.Range("E20").Value = '  +10.33%  '
$ws.Range("D21").NumberFormat = "@"
$ws.Range("D21").Value = '74.46'
$ws.Range("E21").Value = '  +3.26%  '
$ws.Range("E22").Value = '  +1.66%  '
$ws.Range("D23").NumberFormat = "@"
$ws.Range("D23").Value = '265.60'
$ws.Range("E23").Value = '  +9.43%  '
$ws.Range("E24").Value = '  +0.09%  '
$ws.Range("D25").NumberFormat = "@"
$ws.Range("D25").Value = '2.52'
$ws.Range("E25").Value = '  +7.04%  '
$ws.Range("E26").Value = '  +2.08%  '
$ws.Range("D27").NumberFormat = "@"
$ws.Range("D27").Value = '10.12'
$ws.Range("E27").Value = '  +4.98%  '
$ws.Range("D28").NumberFormat = "@"
$ws.Range("D28").Value = '172.98'
$ws.Range("E28").Value = '  +2.19%  '
$ws.Range("D29").NumberFormat = "@"
$ws.Range("D29").Value = '21.20'
$ws.Range("E29").Value = '  +7.24%  '
$ws.Range("E30").Value = '  -2.17%  '
$ws.Range("E31").Value = '  -0.97%  '
$ws.Range("E32").Value = '  +8.16%  '
$ws.Range("D33").NumberFormat = "@"
$ws.Range("D33").Value = '0.125'
$ws.Range("E33").Value = '  +3.05%  '
$ws.Range("E34").Value = '  +6.36%  '
$ws.Range("E35").Value = '  +1.26%  '
$ws.Range("D36").NumberFormat = "@"
$ws.Range("D36").Value = '4.79'
$ws.Range("D37").NumberFormat = "@"
$ws.Range("D37").Value = '6.79'
$ws.Range("E37").Value = '  +7.07%  '
$ws.Range("D38").NumberFormat = "@"
$ws.Range("D38").Value = '3.81'
$ws.Range("E38").Value = '  +7.70%  '
$ws.Range("E39").Value = '  +0.19%  '
$ws.Range("E40").Value = '  +4.67%  '
$ws.Range("E41").Value = '  +0.20%  '
$ws.Range("B42").Value = 'FTXToken'
$ws.Range("C42").Value = 'https://coinranking.com/coin/NfeOYfNcl+ftxtoken-ftt'
$ws.Range("D42").NumberFormat = "@"
$ws.Range("D42").Value = '4.60'
$ws.Range("E42").Value = '  +4.76%  '
$ws.Range("B43").Value = 'FraxShare'
$ws.Range("C43").Value = 'https://coinranking.com/coin/3nNpuxHJ8+fraxshare-fxs'
$ws.Range("D43").NumberFormat = "@"
$ws.Range("D43").Value = '8.44'
$ws.Range("E43").Value = '  -1.39%  '
$ws.Range("D44").NumberFormat = "@"
$ws.Range("D44").Value = '17.53'
$ws.Range("E44").Value = '  +6.86%  '
$ws.Range("E45").Value = '  +2.55%  '
$ws.Range("E46").Value = '  +22.38%  '
$ws.Range("D47").NumberFormat = "@"
$ws.Range("D47").Value = '98.84'
$ws.Range("E47").Value = '  +1.65%  '
$ws.Range("E48").Value = '  +1.43%  '
$ws.Range("B49").Value = 'Maker'
$ws.Range("C49").Value = 'https://coinranking.com/coin/qFakph2rpuMOL+maker-mkr'
$ws.Range("D49").Value = '1.482.79'
$ws.Range("E49").Value = '  +1.09%  '
$ws.Range("B50").Value = 'NEARProtocol'
$ws.Range("C50").Value = 'https://coinranking.com/coin/DCrsaMv68+nearprotocol-near'
$ws.Range("D50").NumberFormat = "@"
$ws.Range("D50").Value = '2.38'
$ws.Range("E50").Value = '  +7.13%  '
$ws.Range("D51").NumberFormat = "@"
$ws.Range("D51").Value = '0.000210'
$ws.Range("E51").Value = '  -12.60%  '
